$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (GFA - Sales), shifting existing
# columns C:J to D:K.
$ws.Columns.Item(3).Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 3).Value = "M_PL"

# New data values for the inserted column (M_PL).
$ws.Cells.Item(2, 3).Value = 1007534436142
$ws.Cells.Item(3, 3).Value = -269766813
$ws.Cells.Item(4, 3).Value = 20228669958
$ws.Cells.Item(5, 3).Value = 344467447608
$ws.Cells.Item(6, 3).Value = 45733381438
